{"js": "// Update the division-problem text cells in the practice worksheet table.\n// Each pair below is (old expression => new expression); every old value\n// is unique in the document, so a case-sensitive whole-document search +\n// replace is unambiguous for each one.\n\nconst replacements = [\n  [\"300\u00f73=\", \"240\u00f75=\"],\n  [\"504\u00f74=\", \"362\u00f73=\"],\n  [\"458\u00f75=\", \"866\u00f74=\"],\n  [\"501\u00f77=\", \"794\u00f74=\"],\n  [\"694\u00f73=\", \"325\u00f75=\"],\n  [\"731\u00f72=\", \"671\u00f75=\"],\n  [\"409\u00f78=\", \"664\u00f77=\"],\n  [\"590\u00f79=\", \"848\u00f73=\"],\n  [\"695\u00f72=\", \"754\u00f73=\"],\n  [\"205\u00f78=\", \"909\u00f75=\"],\n  [\"213\u00f78=\", \"110\u00f72=\"],\n  [\"130\u00f77=\", \"926\u00f77=\"],\n  [\"180\u00f75=\", \"603\u00f77=\"],\n  [\"755\u00f76=\", \"622\u00f78=\"],\n  [\"763\u00f79=\", \"902\u00f76=\"],\n  [\"129\u00f72=\", \"128\u00f76=\"],\n  [\"265\u00f75=\", \"790\u00f75=\"],\n  [\"932\u00f73=\", \"610\u00f77=\"],\n  [\"151\u00f77=\", \"589\u00f79=\"],\n  [\"256\u00f72=\", \"829\u00f77=\"],\n  [\"518\u00f78=\", \"666\u00f76=\"],\n  [\"835\u00f78=\", \"296\u00f74=\"],\n  [\"385\u00f73=\", \"756\u00f77=\"],\n  [\"782\u00f77=\", \"921\u00f78=\"],\n  [\"870\u00f73=\", \"379\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division-problem text cells in the practice worksheet table.\n# Each pair below is (old expression => new expression); every old value\n# is unique in the document, so a plain whole-word Find/Replace over the\n# whole document body is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"300\u00f73=\", \"240\u00f75=\"),\n    @(\"504\u00f74=\", \"362\u00f73=\"),\n    @(\"458\u00f75=\", \"866\u00f74=\"),\n    @(\"501\u00f77=\", \"794\u00f74=\"),\n    @(\"694\u00f73=\", \"325\u00f75=\"),\n    @(\"731\u00f72=\", \"671\u00f75=\"),\n    @(\"409\u00f78=\", \"664\u00f77=\"),\n    @(\"590\u00f79=\", \"848\u00f73=\"),\n    @(\"695\u00f72=\", \"754\u00f73=\"),\n    @(\"205\u00f78=\", \"909\u00f75=\"),\n    @(\"213\u00f78=\", \"110\u00f72=\"),\n    @(\"130\u00f77=\", \"926\u00f77=\"),\n    @(\"180\u00f75=\", \"603\u00f77=\"),\n    @(\"755\u00f76=\", \"622\u00f78=\"),\n    @(\"763\u00f79=\", \"902\u00f76=\"),\n    @(\"129\u00f72=\", \"128\u00f76=\"),\n    @(\"265\u00f75=\", \"790\u00f75=\"),\n    @(\"932\u00f73=\", \"610\u00f77=\"),\n    @(\"151\u00f77=\", \"589\u00f79=\"),\n    @(\"256\u00f72=\", \"829\u00f77=\"),\n    @(\"518\u00f78=\", \"666\u00f76=\"),\n    @(\"835\u00f78=\", \"296\u00f74=\"),\n    @(\"385\u00f73=\", \"756\u00f77=\"),\n    @(\"782\u00f77=\", \"921\u00f78=\"),\n    @(\"870\u00f73=\", \"379\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n}\n"}
